$d = $word.ActiveDocument

# Update Start time
$d.Content.Find.Execute("Start time: 2017-12-27 19:41:36", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start time: 2018-01-31 13:29:32", 2)

# Update End time
$d.Content.Find.Execute("End time: 2017-12-27 19:46:00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "End time: 2018-01-31 13:34:02", 2)

# Update Duration
$d.Content.Find.Execute("Duration: 4.40 mins", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Duration: 4.49 mins", 2)
